$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Customer_Reviews") values were regenerated (scaled ~0.8012x from the
# originals, e.g. to remove multicollinearity with Marketing_Spend) - write the
# exact new values for rows 2 through 51.
$newC = @(
    284.4183153784754,
    643.3462176025796,
    990.2564445290018,
    1286.692435205159,
    1648.0238724888,
    1873.154989732044,
    2259.322955964227,
    2536.530666164093,
    2848.990223903827,
    3181.479240473031,
    1051.145999370591,
    2706.380477038,
    2490.062321679723,
    947.7939918105252,
    2812.937198010781,
    2718.398152335682,
    2459.617544258928,
    1126.456764569399,
    1036.724789013372,
    2294.574803504095,
    1055.953069489664,
    2151.965056638267,
    721.8616962141024,
    2435.582193663564,
    2288.966555031843,
    2943.529269578926,
    2533.325952751378,
    1442.922214075026,
    1971.699927173037,
    1226.604058716749,
    2040.601265546414,
    2051.817762490917,
    1970.097570466679,
    2590.209615827073,
    1823.481931834958,
    3001.2141110078,
    2875.429109558728,
    1203.369886474564,
    799.5759964724464,
    2138.345024634228,
    967.0222722868165,
    1752.177058402044,
    2711.988725510252,
    1616.777916714826,
    2267.334739496016,
    1286.692435205159,
    794.7689263533736,
    1535.057724690588,
    2027.782411895553,
    2133.537954515155
)

for ($i = 0; $i -lt $newC.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $newC[$i]
}

# Header row (A1:D1) gains a thin box border around each cell.
$ws.Range("A1:D1").Borders.LineStyle = 1
